# feat: Proteger rutas de historial y eliminar test-import
#
# Fills in the "Fecha de OC", "Fecha estimada de culminacion",
# "Estado del proyecto en dias" and "Fecha de culminacion real" columns
# (F:I) of the "Presupuestos presentados" sheet, marks each row as
# invoiced ("si"/"no") in column M, backfills N4, narrows the first
# conditional-formatting rule down to H2 only, and moves the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Presupuestos presentados")

# --- Row 2 -------------------------------------------------------------
$ws.Range("D2").Copy($ws.Range("F2"))
$ws.Range("F2").Value = 45012
$ws.Range("D2").Copy($ws.Range("G2"))
$ws.Range("G2").Value = 45717
$ws.Range("H2").Value = 100
$ws.Range("D2").Copy($ws.Range("I2"))
$ws.Range("I2").Value = 45717
$ws.Range("M2").Value = "si"

# --- Row 3 -------------------------------------------------------------
$ws.Range("D3").Copy($ws.Range("F3"))
$ws.Range("F3").Value = 45176
$ws.Range("D2").Copy($ws.Range("G3"))
$ws.Range("G3").Value = 45718
$ws.Range("H3").Value = 100
$ws.Range("D2").Copy($ws.Range("I3"))
$ws.Range("I3").Value = 45718
$ws.Range("M3").Value = "si"

# --- Row 4 -------------------------------------------------------------
$ws.Range("D4").Copy($ws.Range("F4"))
$ws.Range("F4").Value = 45229
$ws.Range("D2").Copy($ws.Range("G4"))
$ws.Range("G4").Value = 45719
$ws.Range("H4").Value = 100
$ws.Range("D2").Copy($ws.Range("I4"))
$ws.Range("I4").Value = 45719
$ws.Range("M4").Value = "si"
$ws.Range("N4").Value = 1409

# --- Row 5 -------------------------------------------------------------
$ws.Range("D5").Copy($ws.Range("F5"))
$ws.Range("F5").Value = 45243
$ws.Range("D2").Copy($ws.Range("G5"))
$ws.Range("G5").Value = 45720
$ws.Range("H5").Value = 100
$ws.Range("D2").Copy($ws.Range("I5"))
$ws.Range("I5").Value = 45720
$ws.Range("M5").Value = "no"

# --- Row 6 -------------------------------------------------------------
$ws.Range("D6").Copy($ws.Range("F6"))
$ws.Range("F6").Value = 45271
$ws.Range("D2").Copy($ws.Range("G6"))
$ws.Range("G6").Value = 45721
$ws.Range("H6").Value = 100
$ws.Range("D2").Copy($ws.Range("I6"))
$ws.Range("I6").Value = 45721
$ws.Range("M6").Value = "no"

# --- Row 7 -------------------------------------------------------------
$ws.Range("D7").Copy($ws.Range("F7"))
$ws.Range("F7").Value = 44918
$ws.Range("D2").Copy($ws.Range("G7"))
$ws.Range("G7").Value = 45722
$ws.Range("H7").Value = 100
$ws.Range("D2").Copy($ws.Range("I7"))
$ws.Range("I7").Value = 45722
$ws.Range("M7").Value = "si"

# --- Row 8 -------------------------------------------------------------
$ws.Range("D8").Copy($ws.Range("F8"))
$ws.Range("F8").Value = 44918
$ws.Range("D2").Copy($ws.Range("G8"))
$ws.Range("G8").Value = 45723
$ws.Range("H8").Value = 100
$ws.Range("D2").Copy($ws.Range("I8"))
$ws.Range("I8").Value = 45723
$ws.Range("M8").Value = "si"

# --- Conditional formatting: shrink "H2:I2 I3:I8" rule down to "H2" ---
$fcs = $ws.Range("H2").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("H2"))

# --- Selection on the "Presupuestos presentados" sheet -----------------
$ws.Activate()
$ws.Range("L10").Select()
